# Rename the original (and only) worksheet from "Sheet1" to "OpenCV".
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "OpenCV"

# Duplicate the OpenCV sheet (placed right after it) to become the new
# "RGB 15 bins" results sheet - this carries over all the labels, the
# Method/Test-No columns and the column formatting/styles untouched.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "RGB 15 bins"

# The new "RGB 15 bins" sheet hasn't actually been run yet, so its
# result column (D) needs to go back to being blank - only the
# formatting (percentage / highlight styles) should remain.
$ws2.Range("D2:D9").ClearContents()
$ws2.Range("D10").ClearContents()
$ws2.Range("D12:D19").ClearContents()
$ws2.Range("D20").ClearContents()
$ws2.Range("D22:D29").ClearContents()
$ws2.Range("D30").ClearContents()
$ws2.Range("D32:D39").ClearContents()
$ws2.Range("D40").ClearContents()

# Column D on the new sheet is no longer sized to fit specific numbers -
# give it a plain custom width instead of the inherited best-fit width.
$ws2.Columns("D").ColumnWidth = 7.4

# Restore view/selection state on both sheets.
$ws1.Activate()
$ws1.Range("A1:D40").Select()

$ws2.Activate()
$ws2.Range("G41").Select()

# "RGB 15 bins" (the newly added sheet) ends up being the active tab.
$ws2.Activate()
